# #5: cash & deposit done
# Rework the "存款" (deposits) sheet: turn row 1 into a proper header row
# and add the standard metadata columns (property_category, category,
# date, legislator_name, legislator_id, source_file, index) to every
# deposit record, mirroring the layout already used on the other sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")
$landWs = $wb.Worksheets.Item("土地")

# ---------------------------------------------------------------------
# 1. Fix header row 1 (it used to be a stray copy of row 2's data).
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"

# New header cells G1:M1 -- copy header formatting (bold/border/center)
# from B1 first, then fill in the labels.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("G1:M1").PasteSpecial(-4122) | Out-Null

$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# ---------------------------------------------------------------------
# 2. Re-point the existing per-row bank / deposit_type / currency cells
#    to their correct values.
# ---------------------------------------------------------------------
$banks   = @("臺灣銀行新營分行", "臺灣銀行新營分行", "臺灣土地銀行新營分行", "中華郵政股份有限公司新營分行", "合作金庫商業銀行新店分行", "聯邦商業銀行新店分行")
$dtypes  = @("活期存款", "定期存款", "活期存款", "活期存款", "活期存款", "活期存款")
$indices = @(47, 48, 49, 50, 51, 52)

for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $banks[$i]
    $ws.Cells.Item($r, 3).Value = $dtypes[$i]
    $ws.Cells.Item($r, 4).Value = "新臺幣"
}

# ---------------------------------------------------------------------
# 3. New columns G:M for every data row (2-7).
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 7).Value  = "deposit"    # G: property_category
    $ws.Cells.Item($r, 8).Value  = "normal"     # H: category
    $ws.Cells.Item($r, 10).Value = "陳唐山"       # J: legislator_name
    $ws.Cells.Item($r, 11).Value = 645          # K: legislator_id
    $ws.Cells.Item($r, 12).Value = "tmp38461"   # L: source_file
    $ws.Cells.Item($r, 13).Value = $indices[$i] # M: index (== col A)
}

# Column I (date) needs the literal text "2012-03-28", not an
# auto-converted date serial -- paste the *value* from a cell that
# already holds that exact string as text (土地!K2), one cell at a
# time (PasteSpecial onto a multi-cell range only fills the anchor).
for ($r = 2; $r -le 7; $r++) {
    $landWs.Range("K2").Copy() | Out-Null
    $ws.Cells.Item($r, 9).PasteSpecial(-4163) | Out-Null
}

$wb.Save()
